$wb = $excel.ActiveWorkbook

# Rename the weekly-tax sheets to their monthly equivalents
$wsGeneral = $wb.Worksheets.Item("GeneralTaxRateWeekly")
$wsGeneral.Name = "GeneralTaxRateMonthly"

$wsProcess = $wb.Worksheets.Item("ProcessPayrollForWeeklyTax")
$wsProcess.Name = "ProcessPayrollForMonthlyTax"

# Update the literal text on the "first" sheet that mirrors those sheet names
$wsFirst = $wb.Worksheets.Item("first")
$wsFirst.Range("A3").Value = "GeneralTaxRateMonthly"
$wsFirst.Range("A4").Value = "ProcessPayrollForMonthlyTax"

# Fix the corrupted employee reference: EMP 107 -> EMP 105
$wsGeneral.Range("A2").Value = "DO NOT TOUCH AUTOMATION EMP 105"
$wsProcess.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 105"

$wsTestReports = $wb.Worksheets.Item("TestReports")
$wsTestReports.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 105"

# Restore view/selection state to match the edited workbook
[void]$wsFirst.Range("F5").Select()
[void]$wsProcess.Range("D9").Select()
[void]$wsTestReports.Range("A3:XFD13").Select()

# GeneralTaxRateMonthly ends up as the active tab
[void]$wsGeneral.Range("A6:XFD14").Select()
